# Update "想去人数" (want-to-go count) values in column F across sheets
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1348
$ws.Range("F5").Value = 887
$ws.Range("F7").Value = 1204
$ws.Range("F8").Value = 1513
$ws.Range("F9").Value = 152
$ws.Range("F10").Value = 52
$ws.Range("F11").Value = 716
$ws.Range("F13").Value = 98
$ws.Range("F17").Value = 78
$ws.Range("F18").Value = 5966
$ws.Range("F19").Value = 45
$ws.Range("F20").Value = 5797
$ws.Range("F21").Value = 9782
$ws.Range("F22").Value = 121
$ws.Range("F23").Value = 170
$ws.Range("F24").Value = 176
$ws.Range("F26").Value = 487
$ws.Range("F28").Value = 141
$ws.Range("F29").Value = 4363
$ws.Range("F30").Value = 360

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 141

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1348
$ws.Range("F10").Value = 1204
$ws.Range("F12").Value = 1513
$ws.Range("F14").Value = 152
$ws.Range("F15").Value = 717
$ws.Range("F18").Value = 98
$ws.Range("F23").Value = 78
$ws.Range("F24").Value = 5966
$ws.Range("F25").Value = 45
$ws.Range("F26").Value = 5797
$ws.Range("F27").Value = 9782
$ws.Range("F29").Value = 121
$ws.Range("F30").Value = 170
$ws.Range("F31").Value = 176
$ws.Range("F34").Value = 487
$ws.Range("F38").Value = 141
$ws.Range("F39").Value = 4363
$ws.Range("F46").Value = 360

